$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 104
$ws.Range("C3").Value = 72
$ws.Range("D3").Value = 117
$ws.Range("F3").Value = 115
$ws.Range("H3").Value = 128
$ws.Range("I3").Value = 175
$ws.Range("J3").Value = 197
$ws.Range("L3").Value = 219
$ws.Range("I6").Value = 21
$ws.Range("B9").Value = 336
$ws.Range("C9").Value = 423
$ws.Range("E9").Value = 395
$ws.Range("F9").Value = 456
$ws.Range("H9").Value = 400
$ws.Range("I9").Value = 447
$ws.Range("K9").Value = 443
$ws.Range("L9").Value = 394
$ws.Range("B10").Value = 1165
$ws.Range("C10").Value = 1374
$ws.Range("D10").Value = 1577
$ws.Range("E10").Value = 1868
$ws.Range("F10").Value = 1879
$ws.Range("G10").Value = 830
$ws.Range("H10").Value = 519
$ws.Range("I10").Value = 755
$ws.Range("J10").Value = 625
$ws.Range("K10").Value = 603
$ws.Range("L10").Value = 582
$ws.Range("B11").Value = 1623
$ws.Range("C11").Value = 1943
$ws.Range("D11").Value = 2151
$ws.Range("E11").Value = 2460
$ws.Range("F11").Value = 2535
$ws.Range("G11").Value = 1443
$ws.Range("H11").Value = 1162
$ws.Range("I11").Value = 1513
$ws.Range("J11").Value = 1313
$ws.Range("K11").Value = 1395
$ws.Range("L11").Value = 1335

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("F3").Value = 3
$ws.Range("I5").Value = 2
$ws.Range("F7").Value = 11

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 30
$ws.Range("E7").Value = 42
$ws.Range("B9").Value = 54
$ws.Range("C9").Value = 85
$ws.Range("E9").Value = 120
$ws.Range("I9").Value = 85

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K7").Value = 38
$ws.Range("L7").Value = 40
$ws.Range("K9").Value = 89
$ws.Range("L9").Value = 86

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 17
$ws.Range("H3").Value = 15
$ws.Range("I3").Value = 28
$ws.Range("C9").Value = 267
$ws.Range("D9").Value = 457
$ws.Range("E9").Value = 549
$ws.Range("F9").Value = 494
$ws.Range("I9").Value = 175
$ws.Range("J9").Value = 102
$ws.Range("K9").Value = 94
$ws.Range("C10").Value = 315
$ws.Range("D10").Value = 522
$ws.Range("E10").Value = 620
$ws.Range("F10").Value = 561
$ws.Range("H10").Value = 173
$ws.Range("I10").Value = 290
$ws.Range("J10").Value = 207
$ws.Range("K10").Value = 178

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("E6").Value = 14
$ws.Range("H7").Value = 11
$ws.Range("E8").Value = 62
$ws.Range("H8").Value = 23

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("F6").Value = 30
$ws.Range("F7").Value = 54

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("H6").Value = 7
$ws.Range("F7").Value = 39
$ws.Range("H7").Value = 6
$ws.Range("F8").Value = 70
$ws.Range("H8").Value = 20

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("E6").Value = 4
$ws.Range("F7").Value = 20
$ws.Range("I7").Value = 21
$ws.Range("B8").Value = 63
$ws.Range("F8").Value = 123
$ws.Range("I8").Value = 75
$ws.Range("D10").Value = 28
$ws.Range("F20").Value = 11
$ws.Range("D22").Value = 17
$ws.Range("E22").Value = 18
$ws.Range("B28").Value = 89
$ws.Range("E28").Value = 78
$ws.Range("F28").Value = 106
$ws.Range("G28").Value = 77
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 77
$ws.Range("L28").Value = 93
$ws.Range("L29").Value = 30
$ws.Range("B32").Value = 54
$ws.Range("C32").Value = 85
$ws.Range("E32").Value = 120
$ws.Range("I32").Value = 85
$ws.Range("I35").Value = 11
$ws.Range("K36").Value = 89
$ws.Range("L36").Value = 86
$ws.Range("B41").Value = 18
$ws.Range("B47").Value = 49
$ws.Range("E47").Value = 67
$ws.Range("F50").Value = 54
$ws.Range("E51").Value = 6
$ws.Range("D52").Value = 34
$ws.Range("E52").Value = 30
$ws.Range("C53").Value = 315
$ws.Range("D53").Value = 522
$ws.Range("E53").Value = 620
$ws.Range("F53").Value = 561
$ws.Range("H53").Value = 173
$ws.Range("I53").Value = 290
$ws.Range("J53").Value = 207
$ws.Range("K53").Value = 178
$ws.Range("F59").Value = 3
$ws.Range("L61").Value = 2
$ws.Range("I62").Value = 23
$ws.Range("K63").Value = 8
$ws.Range("F65").Value = 70
$ws.Range("H65").Value = 20
$ws.Range("E70").Value = 62
$ws.Range("H70").Value = 23
$ws.Range("H71").Value = 2
$ws.Range("F74").Value = 84
$ws.Range("B77").Value = 73
$ws.Range("F77").Value = 56
$ws.Range("I81").Value = 8
$ws.Range("L83").Value = 14
$ws.Range("J87").Value = 26
$ws.Range("I88").Value = 4
$ws.Range("E90").Value = 10
$ws.Range("B92").Value = 23
$ws.Range("B99").Value = 1623
$ws.Range("C99").Value = 1943
$ws.Range("D99").Value = 2151
$ws.Range("E99").Value = 2460
$ws.Range("F99").Value = 2535
$ws.Range("G99").Value = 1443
$ws.Range("H99").Value = 1162
$ws.Range("I99").Value = 1513
$ws.Range("J99").Value = 1313
$ws.Range("K99").Value = 1395
$ws.Range("L99").Value = 1335

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 11
$ws.Range("B8").Value = 18

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 8

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 7
$ws.Range("J8").Value = 12
$ws.Range("J9").Value = 26

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 7
$ws.Range("F3").Value = 8
$ws.Range("J3").Value = 7
$ws.Range("K7").Value = 35
$ws.Range("L7").Value = 40
$ws.Range("B8").Value = 57
$ws.Range("E8").Value = 50
$ws.Range("G8").Value = 34
$ws.Range("B9").Value = 89
$ws.Range("E9").Value = 78
$ws.Range("F9").Value = 106
$ws.Range("G9").Value = 77
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 77
$ws.Range("L9").Value = 93

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("B7").Value = 47
$ws.Range("E7").Value = 53
$ws.Range("B8").Value = 49
$ws.Range("E8").Value = 67

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L7").Value = 20
$ws.Range("L9").Value = 30

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 6
$ws.Range("L6").Value = 14

$ws = $wb.Worksheets.Item('River North')
$ws.Range("F6").Value = 74
$ws.Range("F7").Value = 84

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("B8").Value = 17
$ws.Range("B9").Value = 23

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 4
$ws.Range("I7").Value = 11

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("E6").Value = 4
$ws.Range("E7").Value = 6

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("D7").Value = 27
$ws.Range("E7").Value = 25
$ws.Range("D8").Value = 34
$ws.Range("E8").Value = 30

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I7").Value = 10
$ws.Range("I8").Value = 23

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("F8").Value = 9
$ws.Range("B9").Value = 48
$ws.Range("B10").Value = 73
$ws.Range("F10").Value = 56

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("F4").Value = 2
$ws.Range("F6").Value = 3

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K4").Value = 2
$ws.Range("D7").Value = 14
$ws.Range("E7").Value = 16
$ws.Range("D8").Value = 17
$ws.Range("E8").Value = 18

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 8

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("F6").Value = 11
$ws.Range("I6").Value = 12
$ws.Range("F7").Value = 20
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("D6").Value = 26
$ws.Range("D7").Value = 28

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("E5").Value = 10
$ws.Range("E6").Value = 10

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("E4").Value = 1
$ws.Range("E6").Value = 4

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("B8").Value = 19
$ws.Range("I8").Value = 29
$ws.Range("F9").Value = 83
$ws.Range("B10").Value = 63
$ws.Range("F10").Value = 123
$ws.Range("I10").Value = 75

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 2

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 1
$ws.Range("I7").Value = 4
